$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (Price / Volume(1h) columns) with the latest snapshot.
# Price values that look numeric (e.g. "316.50", "0.4067") must be forced to
# Text format first, otherwise Excel auto-converts them to numbers and silently
# drops significant trailing/leading zeros (e.g. "316.50" -> 316.5). Values that
# already contain multiple "." separators (e.g. "24.716.41") are never parsed as
# numbers by Excel, so no NumberFormat change is needed for those.

$ws.Range("D2").Value = "24.716.41"
$ws.Range("E2").Value = "  +1.56%  "

$ws.Range("D3").Value = "1.696.22"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.50"
$ws.Range("E5").Value = "  +1.39%  "

$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("E7").Value = "  +0.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4067"
$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.488"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.13"
$ws.Range("E11").Value = "  -2.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08852"
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.251"
$ws.Range("E13").Value = "  -0.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.70"
$ws.Range("E14").Value = "  +2.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.060"
$ws.Range("E15").Value = "  +8.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001322"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").Value = "1.698.11"
$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.12"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07030"
$ws.Range("E19").Value = "  -0.39%  "

$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.037"
$ws.Range("E21").Value = "  +4.18%  "

$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.30"
$ws.Range("E23").Value = "  +0.97%  "

$ws.Range("D24").Value = "24.700.41"
$ws.Range("E24").Value = "  +1.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.272"
$ws.Range("E25").Value = "  +9.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.365"
$ws.Range("E26").Value = "  +2.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.84"
$ws.Range("E27").Value = "  +1.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.59"
$ws.Range("E28").Value = "  +2.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.23"
$ws.Range("E29").Value = "  +1.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.184"
$ws.Range("E30").Value = "  +1.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.549"
$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("D32").Value = "1.884.93"
$ws.Range("E32").Value = "  +0.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.070"
$ws.Range("E33").Value = "  -1.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08603"
$ws.Range("E34").Value = "  -1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.142"
$ws.Range("E35").Value = "  -4.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.40"
$ws.Range("E36").Value = "  +1.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2748"
$ws.Range("E37").Value = "  +1.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.896"
$ws.Range("E38").Value = "  -1.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.45"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09244"
$ws.Range("E40").Value = "  +2.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02724"
$ws.Range("E41").Value = "  -1.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.474"
$ws.Range("E42").Value = "  +0.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7658"
$ws.Range("E43").Value = "  +0.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.02"
$ws.Range("E44").Value = "  +3.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7185"
$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.575"
$ws.Range("E46").Value = "  +5.49%  "

$ws.Range("E47").Value = "  +1.72%  "

$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.328"
$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.40"
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07982"
$ws.Range("E51").Value = "  +0.22%  "
